$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.263.53'
$ws.Cells.Item(2, 5).Value = '  +6.32%  '

$ws.Cells.Item(3, 4).Value = '3.536.10'
$ws.Cells.Item(3, 5).Value = '  +8.63%  '

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.19%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '193.53'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +10.03%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '559.96'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +7.41%  '

$ws.Cells.Item(7, 4).Value = '3.529.87'
$ws.Cells.Item(7, 5).Value = '  +8.54%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.613'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.28%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.17%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.645'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +6.85%  '

$ws.Cells.Item(11, 2).Value = 'Avalanche'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '57.00'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +7.42%  '

$ws.Cells.Item(12, 2).Value = 'Dogecoin'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.152'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +16.06%  '

$ws.Cells.Item(13, 5).Value = '  +8.19%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '9.56'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +6.67%  '

$ws.Cells.Item(15, 4).Value = '4.102.92'
$ws.Cells.Item(15, 5).Value = '  +8.62%  '

$ws.Cells.Item(16, 4).Value = '3.534.45'
$ws.Cells.Item(16, 5).Value = '  +8.49%  '

$ws.Cells.Item(17, 5).Value = '  +5.38%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '18.50'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +7.39%  '

$ws.Cells.Item(19, 4).Value = '67.232.97'
$ws.Cells.Item(19, 5).Value = '  +6.30%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '12.01'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +8.78%  '

$ws.Cells.Item(21, 5).Value = '  +4.84%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '409.50'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +11.56%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '4.00'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +6.90%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '86.06'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +6.84%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +9.17%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '11.29'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.69%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +13.20%  '

$ws.Cells.Item(28, 5).Value = '  +0.57%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '12.13'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +7.72%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '8.96'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +9.38%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '30.72'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +8.29%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '682.66'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.15%  '

$ws.Cells.Item(33, 5).Value = '  +6.84%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '11.89'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +6.76%  '

$ws.Cells.Item(35, 5).Value = '  +7.95%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '60.72'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.45%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '39.36'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +7.76%  '

$ws.Cells.Item(38, 4).Value = '0.0₃0832'
$ws.Cells.Item(38, 5).Value = '  +16.35%  '

$ws.Cells.Item(39, 5).Value = '  +0.09%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.399'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +6.27%  '

$ws.Cells.Item(41, 5).Value = '  +13.24%  '

$ws.Cells.Item(42, 5).Value = '  +21.01%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '3.05'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +18.72%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.04%  '

$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45, 4).Value = '3.033.47'
$ws.Cells.Item(45, 5).Value = '  +5.02%  '

$ws.Cells.Item(46, 2).Value = 'Fetch.AI'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.68'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +8.05%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '3.35'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +11.71%  '

$ws.Cells.Item(48, 5).Value = '  +7.73%  '

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +18.93%  '

$ws.Cells.Item(50, 5).Value = '  +3.52%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.132'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +7.17%  '
